$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.876.21'
$ws.Range("E2").Value = '  -1.88%  '

$ws.Range("D3").Value = '2.277.32'
$ws.Range("E3").Value = '  -3.18%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").Value = "'309.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.17%  '

$ws.Range("D6").Value = "'105.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.41%  '

$ws.Range("D7").Value = "'0.624"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.09%  '

$ws.Range("E8").Value = '  +0.15%  '

$ws.Range("D9").Value = "'0.603"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.09%  '

$ws.Range("D10").Value = "'40.18"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.09%  '

$ws.Range("D11").Value = "'0.0905"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.82%  '

$ws.Range("D12").Value = "'8.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.28%  '

$ws.Range("E13").Value = '  -0.08%  '

$ws.Range("D14").Value = "'0.961"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.53%  '

$ws.Range("D15").Value = "'15.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.39%  '

$ws.Range("D16").Value = '2.624.18'
$ws.Range("E16").Value = '  -3.00%  '

$ws.Range("D17").Value = '2.263.37'
$ws.Range("E17").Value = '  -3.61%  '

$ws.Range("D18").Value = '41.853.32'
$ws.Range("E18").Value = '  -1.88%  '

$ws.Range("D19").Value = "'7.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.72%  '

$ws.Range("E20").Value = '  -2.30%  '

$ws.Range("D21").Value = "'73.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.10%  '

$ws.Range("D22").Value = "'3.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -7.23%  '

$ws.Range("D23").Value = "'255.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.96%  '

$ws.Range("D24").Value = "'2.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.49%  '

$ws.Range("D25").Value = "'9.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -7.61%  '

$ws.Range("E26").Value = '  +0.66%  '

$ws.Range("D27").Value = "'10.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.73%  '

$ws.Range("E28").Value = '  +3.35%  '

$ws.Range("D29").Value = "'22.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.91%  '

$ws.Range("D30").Value = "'165.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.72%  '

$ws.Range("D31").Value = "'35.33"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.08%  '

$ws.Range("D32").Value = "'0.0883"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.57%  '

$ws.Range("D33").Value = "'2.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.03%  '

$ws.Range("E34").Value = '  -5.07%  '

$ws.Range("D35").Value = "'0.129"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.34%  '

$ws.Range("D36").Value = "'0.117"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.76%  '

$ws.Range("D37").Value = "'4.52"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.52%  '

$ws.Range("D38").Value = "'0.0350"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.80%  '

$ws.Range("D39").Value = "'2.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.61%  '

$ws.Range("D40").Value = "'3.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.20%  '

$ws.Range("D41").Value = "'71.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.95%  '

$ws.Range("D42").Value = "'97.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.61%  '

$ws.Range("D43").Value = "'1.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.72%  '

$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").Value = "'0.226"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.27%  '

$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = "'1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.03%  '

$ws.Range("D46").Value = "'12.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.96%  '

$ws.Range("D47").Value = "'111.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -8.05%  '

$ws.Range("D48").Value = "'8.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.98%  '

$ws.Range("D49").Value = "'5.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.57%  '

$ws.Range("D50").Value = "'74.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.26%  '

$ws.Range("D51").Value = '1.553.61'
$ws.Range("E51").Value = '  +0.59%  '
